$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 17 de Agosto de 2020 a las 02:25"

# --- Pure numeric refreshes (country/row mapping unchanged) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 5566579
$ws.Range("C4").Value = 36790
$ws.Range("D4").Value = 2922710
$ws.Range("E4").Value = 2470742
$ws.Range("G4").Value = 521
$ws.Range("H4").Value = 173127

# Row 5: Brasil
$ws.Range("E5").Value = 799862
$ws.Range("G5").Value = 582
$ws.Range("H5").Value = 107879

# Row 39
$ws.Range("B39").Value = 81940
$ws.Range("C39").Value = 1275
$ws.Range("D39").Value = 55001
$ws.Range("E39").Value = 25172
$ws.Range("G39").Value = 21
$ws.Range("H39").Value = 1767

# Row 122
$ws.Range("B122").Value = 3016
$ws.Range("C122").Value = 55
$ws.Range("D122").Value = 2036
$ws.Range("E122").Value = 933
$ws.Range("G122").Value = 5
$ws.Range("H122").Value = 47

# Row 161
$ws.Range("D161").Value = 818
$ws.Range("E161").Value = 52

# --- Reordered countries: Paraguay / Consejo Danes para los Refugiados swap ---

$ws.Range("A87").Value = "Paraguay"
$ws.Range("B87").Value = 9791
$ws.Range("C87").Value = 410
$ws.Range("D87").Value = 6034
$ws.Range("E87").Value = 3619
$ws.Range("G87").Value = 11
$ws.Range("H87").Value = 138

$ws.Range("A88").Value = "Consejo Danes para los Refugiados"
$ws.Range("B88").Value = 9676
$ws.Range("C88").Value = 38
$ws.Range("D88").Value = 8705
$ws.Range("E88").Value = 731
$ws.Range("G88").Value = 1
$ws.Range("H88").Value = 240

# --- Reordered countries: Bahamas moves up ahead of Malta / Burkina Faso / Liberia ---

$ws.Range("A148").Value = "Bahamas"
$ws.Range("B148").Value = 1315
$ws.Range("D148").Value = 189
$ws.Range("E148").Value = 1108
$ws.Range("G148").Value = 1
$ws.Range("H148").Value = 18

$ws.Range("A149").Value = "Malta"
$ws.Range("B149").Value = 1306
$ws.Range("C149").Value = 63
$ws.Range("D149").Value = 749
$ws.Range("E149").Value = 548
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 9

$ws.Range("A150").Value = "Burkina Faso"
$ws.Range("B150").Value = 1267
$ws.Range("C150").Value = 27
$ws.Range("D150").Value = 1013
$ws.Range("E150").Value = 199
$ws.Range("G150").Value = 1
$ws.Range("H150").Value = 55

$ws.Range("A151").Value = "Liberia"
$ws.Range("B151").Value = 1257
$ws.Range("D151").Value = 788
$ws.Range("E151").Value = 387
$ws.Range("H151").Value = 82

# --- Reordered countries: Montserrat / Islas Malvinas swap ---

$ws.Range("A213").Value = "Montserrat"
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
